# Remove the trailing empty "Bibliography"-styled paragraph that sits
# immediately before the final section properties (sectPr), right after
# the page-break paragraph near the end of the front matter.
#
# That paragraph is the very last paragraph in the document body, so
# Word's normal Range.Delete() on just that paragraph's own range is a
# no-op (you cannot remove the final paragraph mark of a section on its
# own). Instead we extend the deletion range to also swallow the
# paragraph mark of the preceding paragraph, which removes the empty
# "Bibliography" paragraph while leaving the previous (page-break)
# paragraph and its own mark intact.

$d = $word.ActiveDocument

function Remove-EmptyBibliographyParagraph($doc) {
    $count = $doc.Paragraphs.Count
    for ($i = $count; $i -ge 2; $i--) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Style.NameLocal -eq "Bibliography" -and $p.Range.Text.Trim().Length -eq 0) {
            $prev = $doc.Paragraphs.Item($i - 1)
            $start = $prev.Range.End - 1
            $end = $p.Range.End
            $doc.Range($start, $end).Delete()
            return $true
        }
    }
    return $false
}

Remove-EmptyBibliographyParagraph $d | Out-Null
